$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Rows($ws, $rows) {
    # Capture the original values of columns B:AD for every row in the cycle.
    $ranges = @()
    $values = @()
    foreach ($r in $rows) {
        $rng = $ws.Range("B" + $r + ":AD" + $r)
        $ranges += $rng
        $v = $rng.Value()
        $values += ,$v
    }

    # new[i] = old[i+1] (with wraparound), i.e. each row takes on the content
    # of the next row in the cycle list.
    $n = $rows.Count
    for ($i = 0; $i -lt $n; $i++) {
        $src = ($i + 1) % $n
        $ranges[$i].Value = $values[$src]
    }
}

# Row groups whose match data (columns B:AD, i.e. everything except the id
# in column A) got reshuffled in the updated odds base.
Rotate-Rows $ws @(2, 4)
Rotate-Rows $ws @(7, 10)
Rotate-Rows $ws @(53, 54)
Rotate-Rows $ws @(89, 90)
Rotate-Rows $ws @(114, 115)
Rotate-Rows $ws @(147, 150, 149, 148)
Rotate-Rows $ws @(152, 153)
Rotate-Rows $ws @(154, 155)
Rotate-Rows $ws @(175, 176, 177)
